$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking row: right-answer mark weight changed from 3 to 5
$ws.Range("B11").Value = 5

# Total row: total correct-answer points changed from 45 to 75
$ws.Range("B12").Value = 75

# Corr/total marks text changed from "42/84" to "75/140"
$ws.Range("E12").Value = "75/140"
